# Product Backlog.xlsx - HotFix update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# C3: Story "Y" -> "N"
$ws.Range("C3").Value = "N"

# Fill in Task rows 4-19 (Estimated Time / Owner / Sprint / Status / Accetable Criteria)
$rows1 = @(4,5,6,7,8,9)
foreach ($r in $rows1) {
    $ws.Cells.Item($r, 6).Value = "2hr"      # F: Estimated Time
    $ws.Cells.Item($r, 7).Value = "Matthew"  # G: Owner
    $ws.Cells.Item($r, 8).Value = 2          # H: Sprint
    $ws.Cells.Item($r, 9).Value = "Done"     # I: Status
}
$ws.Cells.Item(4, 10).Value = "Show the list"
$ws.Cells.Item(5, 10).Value = "Show the list"
$ws.Cells.Item(6, 10).Value = "Show the list"
$ws.Cells.Item(7, 10).Value = "Show the list with Limit"
$ws.Cells.Item(8, 10).Value = "Show the list with Limit"
$ws.Cells.Item(9, 10).Value = "Show the list with Limit"

$rows2 = @(10,11,12,13,14,15,16,17,18,19)
foreach ($r in $rows2) {
    $ws.Cells.Item($r, 6).Value = "2hr"          # F: Estimated Time
    $ws.Cells.Item($r, 7).Value = "Giovanmaria"  # G: Owner
    $ws.Cells.Item($r, 8).Value = 2              # H: Sprint
    $ws.Cells.Item($r, 9).Value = "Done"         # I: Status
}
$ws.Cells.Item(10, 10).Value = "Show the list"
$ws.Cells.Item(11, 10).Value = "Show the list"
$ws.Cells.Item(12, 10).Value = "Show the list"
$ws.Cells.Item(13, 10).Value = "Show the list"
$ws.Cells.Item(14, 10).Value = "Show the list"
$ws.Cells.Item(15, 10).Value = "Show the list with Limit"
$ws.Cells.Item(16, 10).Value = "Show the list with Limit"
$ws.Cells.Item(17, 10).Value = "Show the list with Limit"
$ws.Cells.Item(18, 10).Value = "Show the list with Limit"
$ws.Cells.Item(19, 10).Value = "Show the list with Limit"

# Continue the ID numbering in column A for rows 43-58
$idRows = @{43=41;44=42;45=43;46=44;47=45;48=46;49=47;50=48;51=49;52=50;53=51;54=52;55=53;56=54;57=55;58=56}
foreach ($r in $idRows.Keys) {
    $ws.Cells.Item($r, 1).Value = $idRows[$r]
}

# Update the active selection to match the author's last cursor position
$ws.Range("G16").Select()
